$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Myoc"
$ws.Range("C2").Value = "Fzd10"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2144083333333333
$ws.Range("H2").Value = 0.643225
$ws.Range("I2").Value = 0.008611346839948651
$ws.Range("J2").Value = 0.008611346839948651
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.06694666666666667
$ws.Range("N2").Value = 0.20084
$ws.Range("O2").Value = 0.5868893752684747
$ws.Range("P2").Value = 0.5868893752684747
$ws.Range("Q2").Value = 0.01435392322222222
$ws.Range("R2").Value = 0.129185309
$ws.Range("S2").Value = 0.005053907967117618
$ws.Range("T2").Value = 0.005053907967117618

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Myoc"
$ws.Range("C3").Value = "Fzd10"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2144083333333333
$ws.Range("H3").Value = 0.643225
$ws.Range("I3").Value = 0.008611346839948651
$ws.Range("J3").Value = 0.008611346839948651
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.04712366666666667
$ws.Range("N3").Value = 0.141371
$ws.Range("O3").Value = 0.4131106247315253
$ws.Range("P3").Value = 0.4131106247315252
$ws.Range("Q3").Value = 0.01010370683055556
$ws.Range("R3").Value = 0.090933361475
$ws.Range("S3").Value = 0.003557438872831034
$ws.Range("T3").Value = 0.003557438872831033

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Myoc"
$ws.Range("C4").Value = "Fzd10"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 23.971258
$ws.Range("H4").Value = 71.91377399999999
$ws.Range("I4").Value = 0.9627648963950115
$ws.Range("J4").Value = 0.9627648963950115
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.06694666666666667
$ws.Range("N4").Value = 0.20084
$ws.Range("O4").Value = 0.5868893752684747
$ws.Range("P4").Value = 0.5868893752684747
$ws.Range("Q4").Value = 1.604795818906666
$ws.Range("R4").Value = 14.44316237016
$ws.Range("S4").Value = 0.5650364885756861
$ws.Range("T4").Value = 0.5650364885756861

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Myoc"
$ws.Range("C5").Value = "Fzd10"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 23.971258
$ws.Range("H5").Value = 71.91377399999999
$ws.Range("I5").Value = 0.9627648963950115
$ws.Range("J5").Value = 0.9627648963950115
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.04712366666666667
$ws.Range("N5").Value = 0.141371
$ws.Range("O5").Value = 0.4131106247315253
$ws.Range("P5").Value = 0.4131106247315252
$ws.Range("Q5").Value = 1.129613571572667
$ws.Range("R5").Value = 10.166522144154
$ws.Range("S5").Value = 0.3977284078193254
$ws.Range("T5").Value = 0.3977284078193253

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Myoc"
$ws.Range("C6").Value = "Fzd10"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7126843333333334
$ws.Range("H6").Value = 2.138053
$ws.Range("I6").Value = 0.02862375676503981
$ws.Range("J6").Value = 0.02862375676503981
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.06694666666666667
$ws.Range("N6").Value = 0.20084
$ws.Range("O6").Value = 0.5868893752684747
$ws.Range("P6").Value = 0.5868893752684747
$ws.Range("Q6").Value = 0.04771184050222223
$ws.Range("R6").Value = 0.4294065645200001
$ws.Range("S6").Value = 0.01679897872567099
$ws.Range("T6").Value = 0.01679897872567099

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Myoc"
$ws.Range("C7").Value = "Fzd10"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7126843333333334
$ws.Range("H7").Value = 2.138053
$ws.Range("I7").Value = 0.02862375676503981
$ws.Range("J7").Value = 0.02862375676503981
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.04712366666666667
$ws.Range("N7").Value = 0.141371
$ws.Range("O7").Value = 0.4131106247315253
$ws.Range("P7").Value = 0.4131106247315252
$ws.Range("Q7").Value = 0.03358429896255555
$ws.Range("R7").Value = 0.302258690663
$ws.Range("S7").Value = 0.01182477803936882
$ws.Range("T7").Value = 0.01182477803936882

